# Inserts one new data row (for Espinaca / Femacal de La Calera / Coquimbo)
# at row 519, pushing the previously existing rows 519:585 down to 520:586.
# This matches the canonical-XML diff: dimension grows from A1:R585 to
# A1:R586 and a brand-new record appears ahead of the old row 519.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 519:585 down to 520:586, leaving a blank row 519 behind.
$ws.Rows("519").Insert()

# Populate the newly inserted row 519 with the new record's data.
$ws.Range("A519").Value = 3
$ws.Range("B519").Value = "Femacal de La Calera"
$ws.Range("C519").Value = "Coquimbo"
$ws.Range("D519").Value = 45142
$ws.Range("E519").Value = 5
$ws.Range("F519").Value = 100112012
$ws.Range("G519").Value = "Espinaca"
$ws.Range("H519").Value = "Sin especificar"
$ws.Range("I519").Value = "Primera"
$ws.Range("J519").Value = 80
$ws.Range("K519").Value = 4000
$ws.Range("L519").Value = 4000
$ws.Range("M519").Value = 4000
$ws.Range("N519").Value = '$/docena de atados (3 kilos)'
$ws.Range("O519").Value = "Provincia de Quillota"
$ws.Range("P519").Value = 1333
$ws.Range("Q519").Value = 3
$ws.Range("R519").Value = "Hortaliza"
